# LSTM_mv.xlsx update: append the latest rows of predicted/actual FX & oil
# data (2024-09-13 .. 2024-09-19, serials 45335-45344) across the D1_USD,
# D1_EUR, D5_EUR and D1_OIL sheets, and move the "active" tab from D1_OIL
# to D5_EUR (the sheet that was being edited last).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "D1_USD" (sheet1): extend rows 125-128
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("D1_USD")

# Copy the formatting of the last fully-populated row (124) down onto the
# three new data rows (125-127) so dates keep their date format etc.
$ws1.Range("A124:F124").Copy()
$ws1.Range("A125:F127").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A125").Value = 45335
$ws1.Range("B125").Value = 4.0027879999999998
$ws1.Range("C125").Value = 4.0397983000000002
$ws1.Range("D125").Formula = "=B125-C125"
$ws1.Range("E125").Formula = "=IF(D125<0,1,0)"
$ws1.Range("F125").Value = 3.9766759999999999

$ws1.Range("A126").Value = 45336
$ws1.Range("B126").Value = 4.0496049999999997
$ws1.Range("C126").Value = 4.0087447000000003
$ws1.Range("D126").Formula = "=B126-C126"
$ws1.Range("E126").Formula = "=IF(D126<0,1,0)"
$ws1.Range("F126").Value = 3.9856563

$ws1.Range("A127").Value = 45337
$ws1.Range("B127").Value = 4.0427239999999998
$ws1.Range("C127").Value = 4.0576806000000003
$ws1.Range("D127").Formula = "=B127-C127"
$ws1.Range("E127").Formula = "=IF(D127<0,1,0)"
$ws1.Range("F127").Value = 4.0643845000000001

$ws1.Range("C128").Value = 4.0545900000000001
$ws1.Range("F128").Value = 4.0659380000000001

# ---------------------------------------------------------------------
# Sheet "D1_EUR" (sheet3): extend rows 399-402
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("D1_EUR")

$ws3.Range("A398:F398").Copy()
$ws3.Range("A399:F401").PasteSpecial(-4122)  # xlPasteFormats

$ws3.Range("A399").Value = 45335
$ws3.Range("B399").Value = 4.3221699999999998
$ws3.Range("C399").Value = 4.3390446000000003
$ws3.Range("D399").Formula = "=B399-C399"
$ws3.Range("E399").Formula = "=IF(D399<0,1,0)"
$ws3.Range("F399").Value = 4.3574833999999996

$ws3.Range("A400").Value = 45336
$ws3.Range("B400").Value = 4.3365200000000002
$ws3.Range("C400").Value = 4.3137097000000004
$ws3.Range("D400").Formula = "=B400-C400"
$ws3.Range("E400").Formula = "=IF(D400<0,1,0)"
$ws3.Range("F400").Value = 4.3609805000000001

$ws3.Range("A401").Value = 45337
$ws3.Range("B401").Value = 4.3385300000000004
$ws3.Range("C401").Value = 4.3698262999999997
$ws3.Range("F401").Value = 4.3698262999999997

$ws3.Range("C402").Value = 4.3523909999999999
$ws3.Range("F402").Value = 4.3590627

# D401 (as typed by the original author) references the *next* row's C
# cell instead of the usual B-C pattern.
$ws3.Range("D401").Formula = "=C402-C401"
$ws3.Range("E401").Formula = "=IF(D401<0,1,0)"

# ---------------------------------------------------------------------
# Sheet "D5_EUR" (sheet5): rows 92-94 get real predictions, rows 95-99
# are appended with "Nan" placeholders (predictions not available yet)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("D5_EUR")

$ws5.Range("B92").Value = 4.3221699999999998
$ws5.Range("B93").Value = 4.3365200000000002
$ws5.Range("B94").Value = 4.3385300000000004

$ws5.Range("A94:C94").Copy()
$ws5.Range("A95:C99").PasteSpecial(-4122)  # xlPasteFormats

$ws5.Range("A95").Value = 45338
$ws5.Range("B95").Value = "Nan"
$ws5.Range("C95").Value = 4.3211874999999997

$ws5.Range("A96").Value = 45341
$ws5.Range("B96").Value = "Nan"
$ws5.Range("C96").Value = 4.3091654999999998

$ws5.Range("A97").Value = 45342
$ws5.Range("B97").Value = "Nan"
$ws5.Range("C97").Value = 4.2970160000000002

$ws5.Range("A98").Value = 45343
$ws5.Range("B98").Value = "Nan"
$ws5.Range("C98").Value = 4.2957830000000001

$ws5.Range("A99").Value = 45344
$ws5.Range("B99").Value = "Nan"
$ws5.Range("C99").Value = 4.2895380000000003

# ---------------------------------------------------------------------
# Sheet "D1_OIL" (sheet6): extend rows 71-74
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("D1_OIL")

$ws6.Range("A70:F70").Copy()
$ws6.Range("A71:F73").PasteSpecial(-4122)  # xlPasteFormats

$ws6.Range("A71").Value = 45335
$ws6.Range("B71").Value = 77.870002999999997
$ws6.Range("C71").Value = 67.919700000000006
$ws6.Range("D71").Formula = "=B71-C71"
$ws6.Range("E71").Formula = "=D71/C71"
$ws6.Range("F71").Value = 67.919700000000006

$ws6.Range("A72").Value = 45336
$ws6.Range("B72").Value = 76.639999000000003
$ws6.Range("C72").Value = 66.243799999999993
$ws6.Range("D72").Formula = "=B72-C72"
$ws6.Range("E72").Formula = "=D72/C72"

$ws6.Range("A73").Value = 45337
$ws6.Range("B73").Value = 78.029999000000004
$ws6.Range("C73").Value = 73.223399999999998
$ws6.Range("D73").Formula = "=B73-C73"
$ws6.Range("E73").Formula = "=D73/C73"
$ws6.Range("F73").Value = 73.223399999999998

$ws6.Range("C74").Value = 70.859099999999998

# ---------------------------------------------------------------------
# Active sheet moved from D1_OIL to D5_EUR, selection parked on E95
# ---------------------------------------------------------------------
$ws5.Activate()
$ws5.Range("E95").Select()

Write-Output "done"
